$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns so Excel does not
# auto-convert numeric-looking strings (e.g. "530.65", "1.00") into numbers.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = '58.597.52'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '3.149.42'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '530.65'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").Value = '139.43'
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.534'
$ws.Range("E8").Value = '  +14.28%  '
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("E10").Value = '  +4.59%  '
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("E12").Value = '  +2.95%  '
$ws.Range("D13").Value = '3.693.84'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").Value = '25.78'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +3.71%  '
$ws.Range("D16").Value = '58.647.01'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D18").Value = '3.140.23'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '12.99'
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").Value = '8.13'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").Value = '371.01'
$ws.Range("E21").Value = '  +2.84%  '
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '0.523'
$ws.Range("E24").Value = '  +3.38%  '
$ws.Range("D25").Value = '69.67'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").Value = '8.20'
$ws.Range("E28").Value = '  +11.76%  '
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").Value = '22.06'
$ws.Range("E31").Value = '  +2.39%  '
$ws.Range("D32").Value = '6.08'
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("D33").Value = '5.15'
$ws.Range("E33").Value = '  +2.62%  '
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("D35").Value = '158.67'
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("D36").Value = '6.26'
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("E37").Value = '  +4.17%  '
$ws.Range("D38").Value = '25.06'
$ws.Range("E38").Value = '  -3.57%  '
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("E40").Value = '  +2.14%  '
$ws.Range("D41").Value = '2.632.67'
$ws.Range("E41").Value = '  +4.85%  '
$ws.Range("D42").Value = '4.24'
$ws.Range("E42").Value = '  +5.62%  '
$ws.Range("D43").Value = '38.98'
$ws.Range("E43").Value = '  +4.06%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.708'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0285'
$ws.Range("E45").Value = '  +5.87%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '3.193.25'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  +12.73%  '
$ws.Range("E49").Value = '  +2.14%  '
$ws.Range("D50").Value = '0.978'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '20.25'
$ws.Range("E51").Value = '  +2.08%  '

# Restore default cell style (remove the temporary text format) so the
# untouched style metadata matches the original workbook.
$textRange.Style = "Normal"

